$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Find the last used row in column A (the "Serie" date column) and append
# the two new daily records right after it, matching the existing layout
# (columns A:E only; F/G are header-only columns in this sheet).
$lastRow = $ws.Cells($ws.Rows.Count, 1).End(-4162).Row

$newRows = @(
    @("04-08-2021", 50, 207, 50, 774),
    @("05-08-2021", 40, 158, 40, 776)
)

$r = $lastRow
foreach ($row in $newRows) {
    $r = $r + 1

    # Column A holds dd-mm-yyyy date labels stored as literal text (shared
    # strings) in this sheet, not as real dates. Assigning the literal
    # string straight to .Value makes the host "smart-type" it into a real
    # date serial (and tags the cell with a new Text number-format style),
    # which doesn't match the source rows (plain shared-string text, no
    # cell style at all). Routing the text through a TEXT() formula avoids
    # the date auto-detection entirely, then Copy/PasteSpecial(values)
    # collapses the formula down to its plain string result in place -
    # leaving a cell that's indistinguishable from the existing rows.
    $dateCell = $ws.Cells.Item($r, 1)
    $dateCell.Formula = '=TEXT("' + $row[0] + '","@")'
    $dateCell.Copy()
    $dateCell.PasteSpecial(-4163)  # xlPasteValues

    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
}

$excel.CutCopyMode = 0
